$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("UI open and close", "Plays whenever the player opens or closes the pause menu", "UIMenuOpen&Close.wav", "Soundly", "Moderate", "N/A"),
    @("UI open and close (2)", "Another SFX for UI open and close", "UIMenuOpenorClose.wav", "Soundly", "Moderate", "N/A"),
    @("Victory sound", "Plays a victory sound bite whenever the player completes the game", "VictorySFX1-4.wav", "Soundly", "High", "N/A"),
    @("Opening and closing doors", "Plays a sound for when the player opens a door to a new section", "DoorOpen&Close.wav", "Soundly", "Moderate", "N/A"),
    @("Opening and closing doors (2)", "Another SFX for Opening and closing doors", "Door(Short).wav", "Soundly", "Moderate", "N/A")
)

$row = 14
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

$ws.Range("C19").Select()
